# Update the "Latest_update" timestamps (column J) for a few projects.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("J6").Value = "2022-11-16T20:26:44Z"
$ws.Range("J11").Value = "2022-11-16T20:46:22Z"
$ws.Range("J19").Value = "2022-11-16T17:47:53Z"

# Update the "Project_test" flags (column K). These are stored as plain text
# ("True"/"False"), not native Excel booleans, so a direct .Value assignment
# (which Excel auto-coerces into a Boolean for the literal strings
# True/False) would change the underlying cell type. Instead, build the text
# value via a formula in a scratch cell, then copy/paste-special the
# *value* into the target cells - this keeps the result a plain text cell,
# matching the original file's representation.
$scratch = $ws.Range("Z100")

$scratch.Formula = '="True"'
$scratch.Copy()
$ws.Range("K3").PasteSpecial(-4163)
$ws.Range("K17").PasteSpecial(-4163)

$scratch.Formula = '="False"'
$scratch.Copy()
$ws.Range("K8").PasteSpecial(-4163)
$ws.Range("K23").PasteSpecial(-4163)

$scratch.Clear()
